$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string rows appended to the change log sheet (rows 10-12).
# Columns: A=timestamp B=type C=entity D=entity_id E=parent_entity_id
#          F=variable G=old_value H=new_value I=name

$rows = @(
    @{ Row = 10; EntityId = "ser_pub_loc___variable_1" },
    @{ Row = 11; EntityId = "ser_pub_loc___variable_25" },
    @{ Row = 12; EntityId = "accident_route___variable_4" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1746734601
    $ws.Cells.Item($row, 2).Value = "update"
    $ws.Cells.Item($row, 3).Value = "variable"
    $ws.Cells.Item($row, 4).Value = $r.EntityId
    $ws.Cells.Item($row, 6).Value = "key"

    # Force "new_value" to be written as literal text "1" (not the number 1),
    # matching the rest of the sheet where numeric-looking values are text.
    $newValCell = $ws.Cells.Item($row, 8)
    $newValCell.NumberFormat = "@"
    $newValCell.Value = "1"
    $newValCell.ClearFormats()
}
